# Rename the data sheet and the dashboard sheet to their new display names.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("bike_buyers").Name = "Bike buyers"
$wb.Worksheets.Item("Sheet2").Name = "Bike Sales Dashboard"

# The pivot tables/pivot cache still reference the old sheet name internally;
# point the shared pivot cache's source range at the renamed sheet so the
# cache (and every pivot table built on it) tracks the rename.
$dataSheet = $wb.Worksheets.Item("Sheet1")
$pivotTable = $dataSheet.PivotTables(1)
$pivotTable.SourceData = "'Bike buyers'!A1:N1001"
